$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.82383417500417
$ws.Range("C2").Value = 13.0674232306198
$ws.Range("E2").Value = 17.1883080589143
$ws.Range("F2").Value = 34.6666912861141
$ws.Range("G2").Value = 25.3329766834895
$ws.Range("H2").Value = 13.54091003459897
$ws.Range("I2").Value = 21.71735445182841
$ws.Range("J2").Value = 7.415309053205515
$ws.Range("L2").Value = 13.00004056400038
$ws.Range("M2").Value = 16.19163638291497
$ws.Range("O2").Value = 20.13417930201993

$ws.Range("B3").Value = 14.20982340383126
$ws.Range("C3").Value = 12.95549928984862
$ws.Range("E3").Value = 17.24274035082043
$ws.Range("F3").Value = 34.77600757688738
$ws.Range("G3").Value = 25.4748207047706
$ws.Range("H3").Value = 13.59965834451468
$ws.Range("I3").Value = 21.85529217857173
$ws.Range("J3").Value = 7.3946921942755
$ws.Range("L3").Value = 12.98850142457487
$ws.Range("M3").Value = 16.03416379300561
$ws.Range("O3").Value = 20.23852890585799

$ws.Range("B4").Value = 13.81877283369786
$ws.Range("C4").Value = 12.88671418141949
$ws.Range("E4").Value = 17.27836802447111
$ws.Range("F4").Value = 34.85188497816137
$ws.Range("G4").Value = 25.57191059466914
$ws.Range("H4").Value = 13.63811774880043
$ws.Range("I4").Value = 21.94463690399509
$ws.Range("J4").Value = 7.381998140936982
$ws.Range("L4").Value = 12.98284871180751
$ws.Range("M4").Value = 15.93785776714597
$ws.Range("O4").Value = 20.30753339061352

$ws.Range("B5").Value = 13.65608274983104
$ws.Range("C5").Value = 12.85868688166846
$ws.Range("E5").Value = 17.29344211352719
$ws.Range("F5").Value = 34.88500206687203
$ws.Range("G5").Value = 25.61397364989348
$ws.Range("H5").Value = 13.65439089710027
$ws.Range("I5").Value = 21.98221672046068
$ws.Range("J5").Value = 7.376817644665323
$ws.Range("L5").Value = 12.98090752480893
$ws.Range("M5").Value = 15.89874032451494
$ws.Range("O5").Value = 20.33689161898264

$ws.Range("B6").Value = 13.62887313181504
$ws.Range("C6").Value = 12.85403371955226
$ws.Range("E6").Value = 17.29597873037435
$ws.Range("F6").Value = 34.89063362467274
$ws.Range("G6").Value = 25.62110860360288
$ws.Range("H6").Value = 13.65712932434029
$ws.Range("I6").Value = 21.98852760460812
$ws.Range("J6").Value = 7.375957007544602
$ws.Range("L6").Value = 12.98060713694329
$ws.Range("M6").Value = 15.8922535604704
$ws.Range("O6").Value = 20.34184125928849

$ws.Range("B7").Value = 13.81659195599828
$ws.Range("C7").Value = 12.88633615734599
$ws.Range("E7").Value = 17.27856906854536
$ws.Range("F7").Value = 34.85232272110345
$ws.Range("G7").Value = 25.57246777790512
$ws.Range("H7").Value = 13.63833478194279
$ws.Range("I7").Value = 21.94513897441315
$ws.Range("J7").Value = 7.381928303902232
$ws.Range("L7").Value = 12.9828210624319
$ws.Range("M7").Value = 15.93732965538669
$ws.Range("O7").Value = 20.30792431485265

$ws.Range("B8").Value = 14.61514933363134
$ws.Range("C8").Value = 13.02885649343443
$ws.Range("E8").Value = 17.20661893372822
$ws.Range("F8").Value = 34.70256286385004
$ws.Range("G8").Value = 25.37980026628693
$ws.Range("H8").Value = 13.56067101538125
$ws.Range("I8").Value = 21.76395122979191
$ws.Range("J8").Value = 7.408208445276649
$ws.Range("L8").Value = 12.99576591476103
$ws.Range("M8").Value = 16.13727952254391
$ws.Range("O8").Value = 20.16913378846549

$ws.Range("B9").Value = 16.06222729934476
$ws.Range("C9").Value = 13.3069964714319
$ws.Range("E9").Value = 17.08299256823812
$ws.Range("F9").Value = 34.47859759706167
$ws.Range("G9").Value = 25.0819872162855
$ws.Range("H9").Value = 13.42730784646696
$ws.Range("I9").Value = 21.4454618544337
$ws.Range("J9").Value = 7.459408350918741
$ws.Range("L9").Value = 13.03241681407768
$ws.Range("M9").Value = 16.53101597380435
$ws.Range("O9").Value = 19.93620916229785

$ws.Range("B10").Value = 17.0445231207391
$ws.Range("C10").Value = 13.50940891985445
$ws.Range("E10").Value = 17.00276276284982
$ws.Range("F10").Value = 34.35684484612225
$ws.Range("G10").Value = 24.91286691311928
$ws.Range("H10").Value = 13.34085518600135
$ws.Range("I10").Value = 21.2338017353697
$ws.Range("J10").Value = 7.49674655205739
$ws.Range("L10").Value = 13.06607323540924
$ws.Range("M10").Value = 16.81941570191905
$ws.Range("O10").Value = 19.78913880649086

$ws.Range("B11").Value = 17.47245730324479
$ws.Range("C11").Value = 13.60083549219214
$ws.Range("E11").Value = 16.96855498415846
$ws.Range("F11").Value = 34.31080329802973
$ws.Range("G11").Value = 24.8469023089423
$ws.Range("H11").Value = 13.30402651538476
$ws.Range("I11").Value = 21.1423387844516
$ws.Range("J11").Value = 7.513656249177044
$ws.Range("L11").Value = 13.08281308478401
$ws.Range("M11").Value = 16.95003043940758
$ws.Range("O11").Value = 19.7274867519988

$ws.Range("B12").Value = 17.63169201604868
$ws.Range("C12").Value = 13.63534303463623
$ws.Range("E12").Value = 16.95592980071955
$ws.Range("F12").Value = 34.29471597898153
$ws.Range("G12").Value = 24.82351490188636
$ws.Range("H12").Value = 13.29043973036514
$ws.Range("I12").Value = 21.10839606832737
$ws.Range("J12").Value = 7.520047210372343
$ws.Range("L12").Value = 13.08935454800556
$ws.Range("M12").Value = 16.99937476560617
$ws.Range("O12").Value = 19.70489842256648

$ws.Range("B13").Value = 17.59752450496198
$ws.Range("C13").Value = 13.62791659035263
$ws.Range("E13").Value = 16.95863425681073
$ws.Range("F13").Value = 34.29812067869467
$ws.Range("G13").Value = 24.82848077003316
$ws.Range("H13").Value = 13.29334990038786
$ws.Range("I13").Value = 21.11567546564743
$ws.Range("J13").Value = 7.518671376854289
$ws.Range("L13").Value = 13.0879367738528
$ws.Range("M13").Value = 16.98875329338142
$ws.Range("O13").Value = 19.70972946661091

$ws.Range("B14").Value = 17.48561453849926
$ws.Range("C14").Value = 13.60367684086917
$ws.Range("E14").Value = 16.96750972139202
$ws.Range("F14").Value = 34.30945275583674
$ws.Range("G14").Value = 24.84494622631838
$ws.Range("H14").Value = 13.30290151835085
$ws.Range("I14").Value = 21.1395324278593
$ws.Range("J14").Value = 7.514182294020872
$ws.Range("L14").Value = 13.08334721583027
$ws.Range("M14").Value = 16.95409255302789
$ws.Range("O14").Value = 19.72561318919417

$ws.Range("B15").Value = 17.41669734523949
$ws.Range("C15").Value = 13.588813897421
$ws.Range("E15").Value = 16.97298896975408
$ws.Range("F15").Value = 34.31656958681275
$ws.Range("G15").Value = 24.85523953611619
$ws.Range("H15").Value = 13.30879897140331
$ws.Range("I15").Value = 21.15423563422839
$ws.Range("J15").Value = 7.511430945022111
$ws.Range("L15").Value = 13.08056224908371
$ws.Range("M15").Value = 16.93284568087357
$ws.Range("O15").Value = 19.73544122332004

$ws.Range("B16").Value = 17.01616754636679
$ws.Range("C16").Value = 13.50341914656512
$ws.Range("E16").Value = 17.00504433227857
$ws.Range("F16").Value = 34.36004214793233
$ws.Range("G16").Value = 24.91739984529816
$ws.Range("H16").Value = 13.34331231944659
$ws.Range("I16").Value = 21.2398759937863
$ws.Range("J16").Value = 7.495639813853806
$ws.Range("L16").Value = 13.06500774808546
$ws.Range("M16").Value = 16.81086519889287
$ws.Range("O16").Value = 19.79327382212323

$ws.Range("B17").Value = 16.76553843057495
$ws.Range("C17").Value = 13.45085096826535
$ws.Range("E17").Value = 17.02529512439468
$ws.Range("F17").Value = 34.38910741699656
$ws.Range("G17").Value = 24.95835322854191
$ws.Range("H17").Value = 13.36512525623044
$ws.Range("I17").Value = 21.29364786664594
$ws.Range("J17").Value = 7.485932032568288
$ws.Range("L17").Value = 13.05582950462814
$ws.Range("M17").Value = 16.73586179431123
$ws.Range("O17").Value = 19.83009919287086

$ws.Range("B18").Value = 16.61960726118326
$ws.Range("C18").Value = 13.42055496428433
$ws.Range("E18").Value = 17.03715834747496
$ws.Range("F18").Value = 34.40670432724318
$ws.Range("G18").Value = 24.9829401570527
$ws.Range("H18").Value = 13.37790670003685
$ws.Range("I18").Value = 21.32502998609387
$ws.Range("J18").Value = 7.480341344606763
$ws.Range("L18").Value = 13.05068515564438
$ws.Range("M18").Value = 16.69266921130548
$ws.Range("O18").Value = 19.85177424598132

$ws.Range("B19").Value = 16.56989557738527
$ws.Range("C19").Value = 13.41028753574058
$ws.Range("E19").Value = 17.04121206287093
$ws.Range("F19").Value = 34.41281323828034
$ws.Range("G19").Value = 24.99144159532372
$ws.Range("H19").Value = 13.38227467475502
$ws.Range("I19").Value = 21.33573343599222
$ws.Range("J19").Value = 7.478447271907053
$ws.Range("L19").Value = 13.04896660016723
$ws.Range("M19").Value = 16.6780369319723
$ws.Range("O19").Value = 19.85919781690077

$ws.Range("B20").Value = 16.79240286045063
$ws.Range("C20").Value = 13.45645331446594
$ws.Range("E20").Value = 17.02311709376464
$ws.Range("F20").Value = 34.38592233158582
$ws.Range("G20").Value = 24.95388679363054
$ws.Range("H20").Value = 13.36277888710485
$ws.Range("I20").Value = 21.2878767841244
$ws.Range("J20").Value = 7.486966180728796
$ws.Range("L20").Value = 13.0567926218425
$ws.Range("M20").Value = 16.74385172231552
$ws.Range("O20").Value = 19.82612791371298

$ws.Range("B21").Value = 17.51856227848384
$ws.Range("C21").Value = 13.61079988680355
$ws.Range("E21").Value = 16.96489387015349
$ws.Range("F21").Value = 34.30608764760341
$ws.Range("G21").Value = 24.84006660420336
$ws.Range("H21").Value = 13.30008622104308
$ws.Range("I21").Value = 21.13250627903607
$ws.Range("J21").Value = 7.515501194294829
$ws.Range("L21").Value = 13.08468981134088
$ws.Range("M21").Value = 16.96427668058937
$ws.Range("O21").Value = 19.72092716079554

$ws.Range("B22").Value = 17.97670942423193
$ws.Range("C22").Value = 13.71100242911386
$ws.Range("E22").Value = 16.92875644740832
$ws.Range("F22").Value = 34.26176744897496
$ws.Range("G22").Value = 24.77496397655776
$ws.Range("H22").Value = 13.26120803988472
$ws.Range("I22").Value = 21.03499756229473
$ws.Range("J22").Value = 7.534077645792538
$ws.Range("L22").Value = 13.10410070936297
$ws.Range("M22").Value = 17.10764225610378
$ws.Range("O22").Value = 19.65659174495469

$ws.Range("B23").Value = 17.73371945927919
$ws.Range("C23").Value = 13.65759044077604
$ws.Range("E23").Value = 16.94786866542674
$ws.Range("F23").Value = 34.28470195709956
$ws.Range("G23").Value = 24.80885619936807
$ws.Range("H23").Value = 13.28176632637694
$ws.Range("I23").Value = 21.08667095964923
$ws.Range("J23").Value = 7.524170205957638
$ws.Range("L23").Value = 13.09363397696724
$ws.Range("M23").Value = 17.03119977167112
$ws.Range("O23").Value = 19.69052343754823

$ws.Range("B24").Value = 16.78026317917763
$ws.Range("C24").Value = 13.45392072228453
$ws.Range("E24").Value = 17.02410109368179
$ws.Range("F24").Value = 34.38735954638747
$ws.Range("G24").Value = 24.95590282333536
$ws.Range("H24").Value = 13.36383893021016
$ws.Range("I24").Value = 21.29048443279024
$ws.Range("J24").Value = 7.486498671991779
$ws.Range("L24").Value = 13.05635678407794
$ws.Range("M24").Value = 16.74023969798138
$ws.Range("O24").Value = 19.82792176037239

$ws.Range("B25").Value = 15.68442937415515
$ws.Range("C25").Value = 13.23201172401517
$ws.Range("E25").Value = 17.11457198307792
$ws.Range("F25").Value = 34.531690639674
$ws.Range("G25").Value = 25.15389258374428
$ws.Range("H25").Value = 13.46136051961372
$ws.Range("I25").Value = 21.52769172444755
$ws.Range("J25").Value = 7.445600883316331
$ws.Range("L25").Value = 13.02130944025232
$ws.Range("M25").Value = 16.42451879598883
$ws.Range("O25").Value = 19.73544122332004
